# MuleGuard consolidated report update — "updates till config standardisation"
#
# A new API scan ("edgetest_config") is folded into the summary table and the
# pass/fail counters for the table are refreshed; the trailing TOTAL row is
# pushed down one row and recomputed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- make room: push the TOTAL row down to row 12 (carries its style) ------
$ws.Range("A11:E11").Copy($ws.Range("A12:E12"))

# --- turn (old) row 11 into a normal data row, matching row 10's style -----
$ws.Range("A10:F10").Copy($ws.Range("A11:F11"))

# --- rewrite the data rows (API name, counters, status, report path) -------
$ws.Range("A2").Value = "customerOrder"
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "FAIL"
$ws.Range("F2").Value = ".\testData\muleguard-reports\customerOrder\report.html"

$ws.Range("A3").Value = "customerOrderV2"
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "FAIL"
$ws.Range("F3").Value = ".\testData\muleguard-reports\customerOrderV2\report.html"

$ws.Range("A4").Value = "customerOrder_config"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = "FAIL"
$ws.Range("F4").Value = ".\testData\muleguard-reports\customerOrder_config\report.html"

$ws.Range("A5").Value = "customerOrder_V3"
$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "FAIL"
$ws.Range("F5").Value = ".\testData\muleguard-reports\customerOrder_V3\report.html"

$ws.Range("A6").Value = "edgetest_config"
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "FAIL"
$ws.Range("F6").Value = ".\testData\muleguard-reports\edgetest_config\report.html"

$ws.Range("A7").Value = "muleapp1"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = "FAIL"
$ws.Range("F7").Value = ".\testData\muleguard-reports\muleapp1\report.html"

$ws.Range("A8").Value = "muleapp1_config"
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = "FAIL"
$ws.Range("F8").Value = ".\testData\muleguard-reports\muleapp1_config\report.html"

$ws.Range("A9").Value = "muleapp2_config"
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = "FAIL"
$ws.Range("F9").Value = ".\testData\muleguard-reports\muleapp2_config\report.html"

$ws.Range("A10").Value = "rakstestmuleapi"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = "FAIL"
$ws.Range("F10").Value = ".\testData\muleguard-reports\rakstestmuleapi\report.html"

$ws.Range("A11").Value = "rakstestmuleapi_config_20251128"
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = "FAIL"
$ws.Range("F11").Value = ".\testData\muleguard-reports\rakstestmuleapi_config_20251128\report.html"

# --- recompute the TOTAL row (now row 12) -----------------------------------
$ws.Range("A12").Value = "TOTAL"
$ws.Range("B12").Value = 130
$ws.Range("C12").Value = 71
$ws.Range("D12").Value = 59
$ws.Range("E12").Value = "SOME FAILURES"
